# Apply updated cryptocurrency price/volume figures to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "43.187.35"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.322.02"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'302.85"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").Value = "'36.06"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "'17.61"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "2.683.06"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "2.264.56"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "43.091.32"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +7.47%  "
$ws.Range("D20").Value = "'6.24"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "'68.09"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'240.55"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'25.53"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").Value = "'168.52"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'34.28"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'9.21"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "'5.19"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'4.77"
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("D35").Value = "'17.74"
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'1.81"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "1.995.95"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.0290"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("D45").Value = "'10.10"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'17.63"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'76.88"
$ws.Range("E48").Value = "  +9.53%  "
$ws.Range("D49").Value = "'55.02"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  +13.04%  "
$ws.Range("D51").Value = "2.548.19"
$ws.Range("E51").Value = "  +0.75%  "
